$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All changed cells already hold text (inline string) values in the original
# workbook; force text number format first so numeric-looking strings (e.g.
# "13.60", "308.15") are not auto-converted to numbers by the COM layer.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.571.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.244.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +1.71%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.15'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.43'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.574'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.09%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.82'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0807'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.25'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.289.80'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.837'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.60'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '44.417.20'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0959'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.37'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.07'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.80'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '238.50'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.98'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +2.56%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.81'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.45'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.98'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.06'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '152.22'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0801'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.109'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.80'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.42'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.80'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '14.31'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0300'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.757.10'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.193'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '80.55'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '99.44'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '70.13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.88'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.20'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.14'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.87%  '
